$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "97.060.67"

# Row 3
$ws.Range("D3").Value = "3.693.18"
$ws.Range("E3").Value = "  +3.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
Set-TextCell "D5" "240.16"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
Set-TextCell "D6" "1.90"
$ws.Range("E6").Value = "  +9.95%  "

# Row 7
Set-TextCell "D7" "655.07"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
Set-TextCell "D8" "0.426"
$ws.Range("E8").Value = "  -1.82%  "

# Row 9
$ws.Range("E9").Value = "  +3.45%  "

# Row 10
$ws.Range("E10").Value = "  +0.03%  "

# Row 11
$ws.Range("D11").Value = "3.692.21"
$ws.Range("E11").Value = "  +3.10%  "

# Row 12
Set-TextCell "D12" "45.50"
$ws.Range("E12").Value = "  +2.39%  "

# Row 13
$ws.Range("E13").Value = "  +0.88%  "

# Row 14
Set-TextCell "D14" "6.87"
$ws.Range("E14").Value = "  +6.62%  "

# Row 15
$ws.Range("D15").Value = "4.377.08"
$ws.Range("E15").Value = "  +2.97%  "

# Row 16
$ws.Range("E16").Value = "  +2.55%  "

# Row 17
$ws.Range("D17").Value = "96.741.38"

# Row 18
Set-TextCell "D18" "9.07"
$ws.Range("E18").Value = "  +4.09%  "

# Row 19
$ws.Range("D19").Value = "3.696.16"
$ws.Range("E19").Value = "  +3.23%  "

# Row 20
Set-TextCell "D20" "19.30"
$ws.Range("E20").Value = "  +6.73%  "

# Row 21
Set-TextCell "D21" "12.89"
$ws.Range("E21").Value = "  +1.66%  "

# Row 22
$ws.Range("E22").Value = "  -0.31%  "

# Row 23
Set-TextCell "D23" "530.54"
$ws.Range("E23").Value = "  +2.85%  "

# Row 24
Set-TextCell "D24" "3.50"
$ws.Range("E24").Value = "  +0.29%  "

# Row 25
Set-TextCell "D25" "7.12"
$ws.Range("E25").Value = "  +2.58%  "

# Row 26
$ws.Range("E26").Value = "  -1.18%  "

# Row 27
Set-TextCell "D27" "102.58"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28
Set-TextCell "D28" "13.46"
$ws.Range("E28").Value = "  +2.52%  "

# Row 29
$ws.Range("E29").Value = "  -1.97%  "

# Row 30
Set-TextCell "D30" "12.53"
$ws.Range("E30").Value = "  +4.08%  "

# Row 31
$ws.Range("E31").Value = "  +1.80%  "

# Row 32
$ws.Range("E32").Value = "  +0.06%  "

# Row 33
$ws.Range("E33").Value = "  +15.17%  "

# Row 34
$ws.Range("E34").Value = "  +0.65%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D35" "32.75"
$ws.Range("E35").Value = "  +2.53%  "

# Row 36
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell "D36" "1.00"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell "D37" "0.609"
$ws.Range("E37").Value = "  +6.86%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D38" "654.25"
$ws.Range("E38").Value = "  +6.37%  "

# Row 39
Set-TextCell "D39" "9.08"
$ws.Range("E39").Value = "  +3.33%  "

# Row 40
Set-TextCell "D40" "7.01"
$ws.Range("E40").Value = "  +16.16%  "

# Row 41
$ws.Range("E41").Value = "  +5.60%  "

# Row 42
Set-TextCell "D42" "2.01"
$ws.Range("E42").Value = "  +3.11%  "

# Row 43
$ws.Range("E43").Value = "  +4.38%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D44" "1.00"
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D45" "38.15"
$ws.Range("E45").Value = "  +15.77%  "

# Row 46
Set-TextCell "D46" "0.460"
$ws.Range("E46").Value = "  +9.22%  "

# Row 47
$ws.Range("E47").Value = "  +4.28%  "

# Row 48
Set-TextCell "D48" "2.33"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49
Set-TextCell "D49" "23.64"

# Row 50
$ws.Range("E50").Value = "  +2.57%  "

# Row 51
$ws.Range("E51").Value = "  +3.63%  "
